# EPBDS-12787 No bruteforce implementation.
# The sample test-case cell B8 held the text "= addAll(null, null); "Hello";"
# (a quote-prefixed literal, not a formula, so the leading "=" never gets
# evaluated). Rename the referenced rule from addAll to addAll1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the literal with a leading apostrophe so Excel keeps storing it as
# quote-prefixed text (t="s" + a quotePrefix="1" cell style) instead of
# re-parsing the leading "=" as a formula.
$ws.Range("B8").Value = "'= addAll1(null, null); ""Hello"";"

# Move the active selection to J8, matching where the author left the cursor
# when they saved the workbook.
$ws.Range("J8").Select() | Out-Null
